$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Funktion" for the employee in row 7 from "Arzt" to the
# standardized abbreviation "PDL".
$ws.Range("F7").Value = "PDL"

# Reflect the active selection left behind by the edit.
$ws.Range("K8").Select()
